$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: variable, coeff(B), icc(C, text), SE(D), n(E), name(F), lower(G), upper(H), n_country(I)
$data = @(
    @{ Row=2;  A="(Intercept)";    B=2.007;  C="0.109786837851994"; D=0.1250748451902961;  E=17877; F="Model 1a"; G=1.801251879661963;   H=2.212748120338037;   I=6 },
    @{ Row=3;  A="lrscale";        B=-0.589; C="0.109786837851994"; D=0.02900097998158037; E=17877; F="Model 1a"; G=-0.6367066120696997; H=-0.5412933879303002; I=6 },
    @{ Row=4;  A="age";            B=-0.007; C="0.109786837851994"; D=0.03508013478102524; E=17877; F="Model 1a"; G=-0.06470682171478652; H=0.05070682171478651; I=6 },
    @{ Row=5;  A="educ";           B=0.173;  C="0.109786837851994"; D=0.01913072501512775; E=17877; F="Model 1a"; G=0.1415299573501148;   H=0.2044700426498851;  I=6 },
    @{ Row=6;  A="polint";         B=-0.355; C="0.109786837851994"; D=0.02179640376058222; E=17877; F="Model 1a"; G=-0.3908550841861577;  H=-0.3191449158138422; I=6 },
    @{ Row=7;  A="sexMale";        B=-0.044; C="0.109786837851994"; D=0.01297828236770344; E=17877; F="Model 1a"; G=-0.06534927449487216; H=-0.02265072550512784; I=6 },
    @{ Row=8;  A="surveyevs2008";  B=-0.059; C="0.109786837851994"; D=0.01527798735919149; E=17877; F="Model 1a"; G=-0.08413228920587;    H=-0.03386771079412999; I=6 },
    @{ Row=9;  A="surveywvs1994";  B=-0.15;  C="0.109786837851994"; D=0.02556068738474944; E=17877; F="Model 1a"; G=-0.1920473307479128;  H=-0.1079526692520872;  I=6 },
    @{ Row=10; A="surveywvs1999";  B=0.068;  C="0.109786837851994"; D=0.03369779570038187; E=17877; F="Model 1a"; G=0.01256712607287182;  H=0.1234328739271282;   I=6 },
    @{ Row=11; A="surveywvs2005";  B=0.399;  C="0.109786837851994"; D=0.02192011013000984; E=17877; F="Model 1a"; G=0.3629414188361338;   H=0.4350585811638662;   I=6 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    # Column C must be stored as text (matches source which uses inline string of the numeric text)
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
    $ws.Cells.Item($r, 9).Value = $item.I
}
